$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1003.75
$ws.Range("I43").Value = 575
$ws.Range("J43").Value = 1432.5
$ws.Range("K43").Value = 575
$ws.Range("L43").Value = 1432.5
$ws.Range("M43").Value = -506
$ws.Range("N43").Value = -1570.5

$ws.Range("H62").Value = 3806.476
$ws.Range("I62").Value = 3244.25
$ws.Range("J62").Value = 4556.1113
$ws.Range("K62").Value = 3244.25
$ws.Range("L62").Value = 4556.1113
$ws.Range("M62").Value = -2620.25
$ws.Range("N62").Value = -5804.1113

$ws.Range("H65").Value = 3806.476
$ws.Range("I65").Value = 3244.25
$ws.Range("J65").Value = 4556.1113
$ws.Range("K65").Value = 16221.25
$ws.Range("L65").Value = 22780.5565
$ws.Range("M65").Value = -13101.25
$ws.Range("N65").Value = -29020.5565

$ws.Range("H88").Value = 7485.5
$ws.Range("I88").Value = 977
$ws.Range("J88").Value = 10739.75
$ws.Range("K88").Value = 977
$ws.Range("L88").Value = 10739.75
$ws.Range("M88").Value = -571
$ws.Range("N88").Value = -11551.75

$ws.Range("H91").Value = 7485.5
$ws.Range("I91").Value = 977
$ws.Range("J91").Value = 10739.75
$ws.Range("K91").Value = 977
$ws.Range("L91").Value = 10739.75
$ws.Range("M91").Value = 427
$ws.Range("N91").Value = -13547.75

$ws.Range("H96").Value = 307.81482
$ws.Range("I96").Value = 165.52942
$ws.Range("J96").Value = 549.7
$ws.Range("K96").Value = 496.58826
$ws.Range("L96").Value = 1649.1
$ws.Range("M96").Value = 876.41174
$ws.Range("N96").Value = -4395.1

$ws.Range("H100").Value = 1580.8
$ws.Range("I100").Value = 1510.4445
$ws.Range("J100").Value = 1686.3334
$ws.Range("K100").Value = 1510.4445
$ws.Range("L100").Value = 1686.3334
$ws.Range("M100").Value = -969.4445000000001
$ws.Range("N100").Value = -2768.3334

$ws.Range("H106").Value = 4645.4614
$ws.Range("I106").Value = 3045
$ws.Range("J106").Value = 8246.5
$ws.Range("K106").Value = 3045
$ws.Range("L106").Value = 8246.5
$ws.Range("M106").Value = -2414
$ws.Range("N106").Value = -9508.5

$ws.Range("H132").Value = 2705.739
$ws.Range("I132").Value = 2343.7896
$ws.Range("J132").Value = 4425
$ws.Range("K132").Value = 7031.3688
$ws.Range("L132").Value = 13275
$ws.Range("M132").Value = -4501.3688
$ws.Range("N132").Value = -18335

$ws.Range("H137").Value = 30584.857
$ws.Range("J137").Value = 74179.07000000001
$ws.Range("L137").Value = 222537.21
$ws.Range("N137").Value = -227637.21

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 448.1613
$ws.Range("I97").Value = 356.20834
$ws.Range("J97").Value = 763.4286
$ws.Range("K97").Value = 356.20834
$ws.Range("L97").Value = 763.4286
$ws.Range("M97").Value = 139.79166
$ws.Range("N97").Value = -1755.4286

$ws.Range("H132").Value = 1643.7106
$ws.Range("I132").Value = 1504.8148
$ws.Range("J132").Value = 1984.6364
$ws.Range("K132").Value = 4514.4444
$ws.Range("L132").Value = 5953.9092
$ws.Range("M132").Value = -1984.4444
$ws.Range("N132").Value = -11013.9092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1761.125
$ws.Range("I134").Value = 1644.375
$ws.Range("K134").Value = 4933.125
$ws.Range("M134").Value = -2398.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1540.102
$ws.Range("I58").Value = 1005.75757
$ws.Range("J58").Value = 2642.1875
$ws.Range("K58").Value = 1005.75757
$ws.Range("L58").Value = 2642.1875
$ws.Range("M58").Value = -802.75757
$ws.Range("N58").Value = -3048.1875

$ws.Range("H132").Value = 1333.963
$ws.Range("I132").Value = 1085.1708
$ws.Range("J132").Value = 2118.6155
$ws.Range("K132").Value = 3255.512400000001
$ws.Range("L132").Value = 6355.8465
$ws.Range("M132").Value = -725.5124000000005
$ws.Range("N132").Value = -11415.8465

$ws.Range("H134").Value = 1975.9783
$ws.Range("I134").Value = 1303.909
$ws.Range("J134").Value = 3682
$ws.Range("K134").Value = 3911.727
$ws.Range("L134").Value = 11046
$ws.Range("M134").Value = -1376.727
$ws.Range("N134").Value = -16116

$ws.Range("H136").Value = 1540.102
$ws.Range("I136").Value = 1005.75757
$ws.Range("J136").Value = 2642.1875
$ws.Range("K136").Value = 3017.27271
$ws.Range("L136").Value = 7926.5625
$ws.Range("M136").Value = -467.2727100000002
$ws.Range("N136").Value = -13026.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 700
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws.Range("H132").Value = 3152.1187
$ws.Range("I132").Value = 3184.4255
$ws.Range("J132").Value = 3025.5833
$ws.Range("K132").Value = 9553.2765
$ws.Range("L132").Value = 9076.749899999999
$ws.Range("M132").Value = -7023.2765
$ws.Range("N132").Value = -14136.7499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 925
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 1300
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 1300
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -1676

$ws.Range("H82").Value = 1843.8096
$ws.Range("I82").Value = 1308.2222
$ws.Range("J82").Value = 2245.5
$ws.Range("K82").Value = 1308.2222
$ws.Range("L82").Value = 2245.5
$ws.Range("M82").Value = -947.2221999999999
$ws.Range("N82").Value = -2967.5

$ws.Range("H85").Value = 1843.8096
$ws.Range("I85").Value = 1308.2222
$ws.Range("J85").Value = 2245.5
$ws.Range("K85").Value = 1308.2222
$ws.Range("L85").Value = 2245.5
$ws.Range("M85").Value = -60.22219999999993
$ws.Range("N85").Value = -4741.5

$ws.Range("H93").Value = 7594
$ws.Range("I93").Value = 13225.75
$ws.Range("J93").Value = 1157.7142
$ws.Range("K93").Value = 13225.75
$ws.Range("L93").Value = 1157.7142
$ws.Range("M93").Value = -11977.75
$ws.Range("N93").Value = -3653.7142

$ws.Range("H122").Value = 3201.087
$ws.Range("I122").Value = 3210.4546
$ws.Range("J122").Value = 3192.5
$ws.Range("K122").Value = 9631.363799999999
$ws.Range("L122").Value = 9577.5
$ws.Range("M122").Value = -7181.363799999999
$ws.Range("N122").Value = -14477.5

$ws.Range("H132").Value = 1399.15
$ws.Range("I132").Value = 1476.4482
$ws.Range("J132").Value = 881.8461
$ws.Range("K132").Value = 4429.3446
$ws.Range("L132").Value = 2645.5383
$ws.Range("M132").Value = -1899.3446
$ws.Range("N132").Value = -7705.5383

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1911
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1911
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1911
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4657

$ws.Range("H100").Value = 698.75
$ws.Range("I100").Value = 550.2857
$ws.Range("K100").Value = 1100.5714
$ws.Range("M100").Value = -559.5714

$ws.Range("H132").Value = 1295.9193
$ws.Range("I132").Value = 773.02563
$ws.Range("J132").Value = 2182.5652
$ws.Range("K132").Value = 2319.07689
$ws.Range("L132").Value = 6547.6956
$ws.Range("M132").Value = 210.9231100000002
$ws.Range("N132").Value = -11607.6956
